$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.002.59"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "'3.172.51"
$ws.Range("E3").Value = "  +4.55%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'216.07"
$ws.Range("E5").Value = "  +2.22%  "

$ws.Range("D6").Value = "'627.03"
$ws.Range("E6").Value = "  +2.35%  "

$ws.Range("D7").Value = "'1.17"
$ws.Range("E7").Value = "  +31.38%  "

$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "'3.170.71"
$ws.Range("E10").Value = "  +4.61%  "

$ws.Range("D11").Value = "'0.761"
$ws.Range("E11").Value = "  +14.40%  "

$ws.Range("E12").Value = "  +7.50%  "

$ws.Range("D13").Value = "'5.72"
$ws.Range("E13").Value = "  +7.65%  "

$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").Value = "'35.02"
$ws.Range("E15").Value = "  +8.39%  "

$ws.Range("D16").Value = "'90.799.52"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").Value = "'3.759.75"
$ws.Range("E17").Value = "  +4.93%  "

$ws.Range("D18").Value = "'3.191.24"
$ws.Range("E18").Value = "  +5.58%  "

$ws.Range("D19").Value = "'3.77"
$ws.Range("E19").Value = "  +13.81%  "

$ws.Range("D20").Value = "'14.69"
$ws.Range("E20").Value = "  +9.45%  "

$ws.Range("D21").Value = "'471.43"
$ws.Range("E21").Value = "  +11.14%  "

$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").Value = "'9.17"
$ws.Range("E23").Value = "  +10.96%  "

$ws.Range("D24").Value = "'5.26"
$ws.Range("E24").Value = "  +4.43%  "

$ws.Range("D25").Value = "'5.93"
$ws.Range("E25").Value = "  +11.29%  "

$ws.Range("D26").Value = "'96.00"
$ws.Range("E26").Value = "  +16.01%  "

$ws.Range("E27").Value = "  +7.37%  "

$ws.Range("D28").Value = "'3.343.51"
$ws.Range("E28").Value = "  +4.98%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("D31").Value = "'9.29"
$ws.Range("E31").Value = "  +9.34%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'28.40"
$ws.Range("E32").Value = "  +24.77%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'0.974"
$ws.Range("E33").Value = "  -20.65%  "

$ws.Range("D34").Value = "'0.194"
$ws.Range("E34").Value = "  +42.30%  "

$ws.Range("D35").Value = "'525.37"
$ws.Range("E35").Value = "  +4.83%  "

$ws.Range("D36").Value = "'1.93"
$ws.Range("E36").Value = "  +6.72%  "

$ws.Range("D37").Value = "'3.63"
$ws.Range("E37").Value = "  -2.63%  "

$ws.Range("D38").Value = "'0.144"
$ws.Range("E38").Value = "  +8.13%  "

$ws.Range("D39").Value = "'6.98"
$ws.Range("E39").Value = "  +5.17%  "

$ws.Range("D40").Value = "'1.31"
$ws.Range("E40").Value = "  +5.46%  "

$ws.Range("D41").Value = "'0.0902"
$ws.Range("E41").Value = "  +29.82%  "

$ws.Range("D42").Value = "'22.23"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("E43").Value = "  +17.23%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'2.00"
$ws.Range("E45").Value = "  +9.54%  "

$ws.Range("D47").Value = "'0.706"
$ws.Range("E47").Value = "  +19.49%  "

$ws.Range("D48").Value = "'150.22"
$ws.Range("E48").Value = "  +4.99%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'4.61"
$ws.Range("E49").Value = "  +9.52%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.37"
$ws.Range("E50").Value = "  +12.28%  "

$ws.Range("D51").Value = "'45.43"
$ws.Range("E51").Value = "  +4.31%  "
